$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.664.52'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '1.888.85'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '237.25'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4838'
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").Value = '0.2856'
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("D9").Value = '0.06541'
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.829.65'
$ws.Range("E10").Value = '  -1.58%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.07446'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '16.54'
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("D13").Value = '5.083'
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("D14").Value = '87.82'
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = '0.6631'
$ws.Range("E15").Value = '  +3.02%  '
$ws.Range("D16").Value = '30.602.85'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = '13.18'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '0.000007595'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = '229.83'
$ws.Range("E20").Value = '  +3.01%  '
$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.260'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.048.66'
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("D24").Value = '6.179'
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").Value = '9.410'
$ws.Range("E25").Value = '  +2.58%  '
$ws.Range("D26").Value = '167.57'
$ws.Range("E26").Value = '  +2.56%  '
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '1.951'
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").Value = '0.1023'
$ws.Range("E29").Value = '  +11.24%  '
$ws.Range("D30").Value = '1.394'
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").Value = '4.329'
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").Value = '4.020'
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").Value = '0.05049'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").Value = '1.201'
$ws.Range("E34").Value = '  +5.26%  '
$ws.Range("D35").Value = '0.7493'
$ws.Range("E35").Value = '  +3.59%  '
$ws.Range("D36").Value = '0.9988'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").Value = '0.01888'
$ws.Range("E38").Value = '  +3.30%  '
$ws.Range("D39").Value = '2.663'
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("D40").Value = '0.9192'
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").Value = '2.061'
$ws.Range("E41").Value = '  +1.28%  '
$ws.Range("D42").Value = '107.13'
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("D43").Value = '0.4263'
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("D45").Value = '5.641'
$ws.Range("E45").Value = '  -4.31%  '
$ws.Range("D46").Value = '7.406'
$ws.Range("E46").Value = '  +2.14%  '
$ws.Range("D47").Value = '64.64'
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").Value = '0.1271'
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").Value = '8.952'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '33.99'
$ws.Range("E51").Value = '  +0.87%  '
